# Update "想去人数" (number of people interested) counts that changed
# between scrapes, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 272
$ws1.Range("F4").Value = 929
$ws1.Range("F6").Value = 48

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 272
$ws4.Range("F5").Value = 929
$ws4.Range("F7").Value = 48
